$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C25").Value = 57900
$ws.Range("E25").Value = "In-house data (Shobhan)"

$ws.Range("C26").Value = 0.00000269
$ws.Range("E26").Value = "In-house data (Shobhan)"

$ws.Range("E26").Select()
